# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handback timestamp cells
# to reflect the new report-generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-24 03:05:39"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-24 03:05:34"
$zhcn.Range("K2").Value = "2016-08-24 03:05:52"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-24 03:05:39"
$dede.Range("K2").Value = "2016-08-24 03:05:59"
